# edit.ps1 -- apply the HW2 report revision described by the commit diff.
#
# Strategy: this runtime's Range.InsertXML reliably replaces a paragraph's
# content (including run-splitting and <w:proofErr/> / <w:bookmarkStart|End/>
# placement) only when the target Range spans exactly one whole paragraph
# (paragraph mark included) AND is obtained via Document.Paragraphs.Item(n).Range
# -- using a Find-derived/Expand(4) range instead corrupts the very last
# paragraph of the body. Plain Find/Replace text substitution cannot create
# new runs or proofErr markers, so each touched paragraph is rewritten whole,
# from literal target OOXML, via the Paragraphs collection. Paragraph count
# is unchanged by this edit, so the indices below are stable.

$d = $word.ActiveDocument

function Set-ParagraphXml($doc, [int]$paraIndex, [string]$expectedLeadingText, [string]$xml) {
    $para = $doc.Paragraphs.Item($paraIndex)
    $rng = $para.Range
    if ($rng.Text.IndexOf($expectedLeadingText) -ne 0) {
        throw "Paragraph $paraIndex does not start with expected text. Found: [$($rng.Text)]"
    }
    $rng.InsertXML($xml)
}

# Paragraph 13: "To implement cross-validation, ..." -- split off "instances,"
# into its own run wrapped in w:proofErr gramStart/gramEnd markers.
$xmlP13 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To implement cross-validation, the training data was first divided into 10 equal parts, each consisting of an equal number of positive and negative </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>instances,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> randomly sorted.  9 parts were used to train, while the last was used to test, rotating through each possible 9:1 combination and averaging the accuracies.  Both a linear model and a 5-degree polynomial model were used for training and testing.</w:t></w:r></w:p>
'@
Set-ParagraphXml $d 13 "To implement cross-validation" $xmlP13

# Paragraph 34: "The variance with the polynomial ... model was more visible"
# -- wraps "over 4% better accuracy ... classifying the test set" in a
# w:proofErr gramStart/gramEnd pair (splitting the first affected run).
$xmlP34 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The variance with the polynomial </w:t></w:r><w:r><w:t xml:space="preserve">cross-validation </w:t></w:r><w:r><w:t xml:space="preserve">model was more visible - </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>over 4% better accuracy</w:t></w:r><w:r><w:t xml:space="preserve"> than that of the full training classifying the test set</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">.  </w:t></w:r><w:r><w:t>Cross-validation actually proved a worthwhile investment of time; t</w:t></w:r><w:r><w:t>his was clearly the more accurate model.  However, all results were still well under 70%, so there is still much room for this to improve.</w:t></w:r></w:p>
'@
Set-ParagraphXml $d 34 "The variance with the polynomial" $xmlP34

# Paragraph 38: "To script the Adaboost algorithm, ..." -- wraps "Adaboost" in
# spellStart/spellEnd proofErr, appends the new sentence about the
# roulette-selection "barrier" array, and the _GoBack bookmark now closes
# this paragraph (relocated from the final CONCLUSION paragraph).
$xmlP38 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To script the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adaboost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> algorithm, I used Perl, from which SVM-Light functions were called.  10 iterations of boosting were implemented, calculating the error, alpha, and weight changes.</w:t></w:r><w:r><w:t xml:space="preserve">  I used a second array to hold each weight’s “barrier” so that each slice was represented as a range somewhere between 0 and 1.  Therefore, when a random number was chosen for the roulette selection, it was more likely to fall within the larger slices, thus focusing more on these.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
Set-ParagraphXml $d 38 "To script the Adaboost algorithm" $xmlP38

# Paragraph 40: "I found that boosting actually produced ..." -- updates the
# accuracy figure to 61.02% and softens "was slightly" to "is still",
# splitting the numbers/phrase out into their own runs.
$xmlP40 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I found that </w:t></w:r><w:r><w:t xml:space="preserve">boosting actually produced lower accuracy rates - about </w:t></w:r><w:r><w:t>61.02</w:t></w:r><w:r><w:t xml:space="preserve">% - than cross-validation.  However, this </w:t></w:r><w:r><w:t>is still</w:t></w:r><w:r><w:t xml:space="preserve"> improved over the original accuracy rates, ranging between 49% and 52%. </w:t></w:r></w:p>
'@
Set-ParagraphXml $d 40 "I found that " $xmlP40

# Paragraph 42: "Increasing the number of iterations ..." -- drops the
# "30, and even 50" clause and the parenthetical, replacing with a single
# 20-iteration ensemble accuracy of 61.8%, split into several runs.
$xmlP42 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Increasing the number of iterations didn’t appear to help much.  20</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">iterations still yielded </w:t></w:r><w:r><w:t>an ensemble accuracy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>of 61.8%.</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r></w:p>
'@
Set-ParagraphXml $d 42 "Increasing the number of iterations" $xmlP42

# Paragraph 46: "Cross-validation had a much greater impact ..." -- wraps
# each "Adaboost" occurrence in spellStart/spellEnd proofErr, and the
# _GoBack bookmark is removed from here (now lives on paragraph 38 instead).
$xmlP46 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Cross-validation had a much greater impact on classification accuracies than </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adaboost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  I believe this has to do with the initial accuracy of classification.  As the base algorithm sits right around 50% (at or below random), it is not surprising that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adaboost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> would be relatively ineffective, as it best</w:t></w:r><w:r><w:t xml:space="preserve"> improves accuracy when augmenting a classification with better-than-random performance.  I might be curious to see if the two techniques could be combined, using the output data from cross-validation to feed into </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Adaboost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.  This could potentially produce better classification.</w:t></w:r></w:p>
'@
Set-ParagraphXml $d 46 "Cross-validation had a much greater impact" $xmlP46

Write-Output "done"
